$d = $word.ActiveDocument

$replacements = @(
    @("77÷2=", "85÷9="),
    @("22÷4=", "70÷8="),
    @("77÷7=", "49÷2="),
    @("53÷3=", "58÷7="),
    @("90÷6=", "28÷8="),
    @("30÷3=", "82÷2="),
    @("81÷9=", "84÷8="),
    @("71÷8=", "46÷6="),
    @("16÷7=", "68÷4="),
    @("62÷9=", "65÷9="),
    @("27÷5=", "86÷4="),
    @("41÷3=", "59÷7="),
    @("21÷6=", "93÷4="),
    @("21÷2=", "22÷2="),
    @("66÷2=", "58÷2="),
    @("84÷3=", "24÷5="),
    @("80÷8=", "85÷3="),
    @("11÷2=", "22÷8="),
    @("35÷5=", "94÷3="),
    @("57÷3=", "35÷6="),
    @("53÷5=", "10÷6="),
    @("63÷3=", "76÷4="),
    @("74÷6=", "42÷2="),
    @("91÷6=", "79÷9="),
    @("12÷6=", "45÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
